{"js": "// Helper: find the first (and expected-only) exact match of `oldText` in the\n// document body and replace its text with `newText`, preserving the run's\n// original formatting (font/size/color) because Range.insertText keeps the\n// formatting of the text it overwrites.\nasync function replaceOnce(context, oldText, newText) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n  const range = results.items[0];\n  range.insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n  return range;\n}\n\n// Helper: find the first exact match of `anchorText`, then insert `newText`\n// immediately after it (used to splice in brand-new sentences that don't\n// exist yet in the original document).\nasync function insertAfterAnchor(context, anchorText, newText) {\n  const results = context.document.body.search(anchorText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Anchor text not found: \" + anchorText);\n  }\n  const anchor = results.items[0];\n  const inserted = anchor.insertText(newText, Word.InsertLocation.after);\n  await context.sync();\n  return inserted;\n}\n\n// 1) Title\nawait replaceOnce(context, \"Quantum Mysteries Unveiled\", \"Government and Public Policy: Unveiling the Interwoven Dynamics\");\n\n// 2) Byline (author name)\nawait replaceOnce(context, \"Isaac Newton\", \"Dr. Henry Morgan\");\n\n// 3) Email line - the original is split across runs (\"isaac\" + \".\" + \"newton@physics\" + \".\" + \"org\"),\n// but together they form the single string \"isaac.newton@physics.org\"; replace it wholesale.\nawait replaceOnce(context, \"isaac.newton@physics.org\", \"hmorgan@moorehighschool.com\");\n\n// 4) First body paragraph - four sentences before the first line break.\nawait replaceOnce(\n  context,\n  \"In the realm of physics, the exploration of quantum mechanics has yielded profound insights, challenging our understanding of the universe\",\n  \"In the realm of human society, nothing is more pivotal than the intricate tapestry of government and public policy\"\n);\nawait replaceOnce(\n  context,\n  \" Quantum mechanics, initially conceived as a framework to explain the enigmatic behavior of subatomic particles, has revolutionized various disciplines, including chemistry, materials science, and even computer science\",\n  \" These two forces, intertwined like threads in a multihued fabric, exert profound influences on the trajectory of our lives\"\n);\nawait replaceOnce(\n  context,\n  \" This essay delves into the captivating world of quantum physics, investigating its history, fundamental principles, and transformative applications\",\n  \" Government shapes societal structure, crafting laws and policies that weave through every aspect of our existence\"\n);\n// New 4th sentence inserted after the period that follows the previous sentence.\nawait insertAfterAnchor(\n  context,\n  \"Government shapes societal structure, crafting laws and policies that weave through every aspect of our existence.\",\n  \" Public policy, in turn, is a vivid reflection of the collective aspirations and values of society, a tapestry woven from the threads of governmental action.\"\n);\n\n// Second block (after first double line-break)\nawait replaceOnce(\n  context,\n  \"This conceptual revolution, sparked by Max Planck's introduction of energy quantization, shattered long-held assumptions about the continuity of energy and revealed the discreet nature of light and matter\",\n  \"These interconnected entities paint a picture of complex interactions\"\n);\nawait replaceOnce(\n  context,\n  \" Albert Einstein's groundbreaking photoelectric effect experiment provided empirical validation for this departure from classical physics\",\n  \" Government draws its legitimacy from the consent of the governed, while public policy becomes a testament to the collective will of the people\"\n);\nawait replaceOnce(\n  context,\n  \" Moreover, the advent of wave-particle duality, epitomized by the double-slit experiment, unveiled the paradoxical behavior of subatomic particles exhibiting both wave-like and particle-like properties\",\n  \" As the government formulates policies that shape regulations, taxation, and public services, it lays the foundation for a just and equitable society\"\n);\nawait insertAfterAnchor(\n  context,\n  \"As the government formulates policies that shape regulations, taxation, and public services, it lays the foundation for a just and equitable society.\",\n  \" These policies impact everything from economic growth and environmental protection to education and healthcare, weaving their way into the very fabric of our daily lives.\"\n);\n\n// Third block (after second double line-break)\nawait replaceOnce(\n  context,\n  \"The enigma of quantum entanglement, wherein particles exhibit a profound interconnectedness regardless of distance, has perplexed scientists and sparked profound debates about the nature of reality\",\n  \"Furthermore, government and public policy are dynamic entities, constantly evolving in response to changing societal needs and global challenges\"\n);\nawait replaceOnce(\n  context,\n  \" This phenomenon, defying classical notions of locality, has profound implications, ranging from cryptography to quantum computing\",\n  \" As the world grapples with issues such as climate change, resource depletion, and globalization, governments must adapt their policies to confront these pressing concerns\"\n);\nawait replaceOnce(\n  context,\n  \" It challenges our understanding of information transfer and raises fundamental questions regarding the relationship between consciousness and the physical world\",\n  \" This delicate dance between government and public policy ensures that society navigates the ever-shifting tides of progress, preserving core values while embracing the imperatives of a changing world\"\n);\n\n// 5) Summary paragraph\nawait replaceOnce(\n  context,\n  \"Quantum mechanics has reshaped our understanding of the universe, providing a conceptual framework for explaining the enigmatic behavior of subatomic particles\",\n  \"In this essay, we explored the intrinsic connection between government and public policy\"\n);\nawait replaceOnce(\n  context,\n  \" The quantization of energy, wave-particle duality, and quantum entanglement challenge classical notions of physics and open up new avenues for exploration across a wide range of fields\",\n  \" We recognized the critical role government plays in shaping society through the policies it enacts\"\n);\n// This sentence used to straddle a page break run (\"... and \" + \"revolutionizing industries\"); now it\n// becomes one self-contained sentence and the following sentence is replaced independently below.\nawait replaceOnce(\n  context,\n  \" The transformative applications of quantum mechanics, from quantum computing to cryptography, hold immense promise for advancing technology and revolutionizing industries\",\n  \" These policies, in turn, are influenced by societal values and aspirations, and form the backbone of a stable and just society\"\n);\nawait replaceOnce(\n  context,\n  \" As we delve deeper into the mysteries of the quantum realm, we unlock new frontiers of knowledge with the potential to shape our future in unimaginable ways\",\n  \" Government and public policy work in tandem, evolving over time to address new Herausforderungen and societal shifts\"\n);\nawait insertAfterAnchor(\n  context,\n  \"Government and public policy work in tandem, evolving over time to address new Herausforderungen and societal shifts.\",\n  \" By fostering this vital partnership, we ensure a government that is responsive to the people it serves and policies that reflect the collective will and best interests of society as a whole.\"\n);\n\n// 6) Trailing empty paragraph added at the end of the document body.\nconst lastParagraph = context.document.body.paragraphs.getLast();\nlastParagraph.insertParagraph(\"\", Word.InsertLocation.after);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Replace the first exact occurrence of $OldText anywhere in the document\n# body with $NewText. Using Find+Range.Text keeps the formatting (font,\n# size, color) of the run(s) being overwritten, same as Word would when a\n# user selects the text and types over it.\nfunction Replace-Text($OldText, $NewText) {\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.MatchCase = $true\n    $rng.Find.MatchWholeWord = $false\n    $rng.Find.Execute($OldText) | Out-Null\n    if (-not $rng.Find.Found) {\n        Write-Output \"NOT FOUND: $OldText\"\n        return\n    }\n    $rng.Text = $NewText\n}\n\n# Find the first exact occurrence of $AnchorText, collapse to its end, and\n# insert $NewText right after it. Used to splice in brand-new sentences\n# that have no counterpart in the original document.\nfunction Insert-After($AnchorText, $NewText) {\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.MatchCase = $true\n    $rng.Find.MatchWholeWord = $false\n    $rng.Find.Execute($AnchorText) | Out-Null\n    if (-not $rng.Find.Found) {\n        Write-Output \"ANCHOR NOT FOUND: $AnchorText\"\n        return\n    }\n    $rng.Collapse(0)\n    $rng.InsertAfter($NewText)\n}\n\n# 1) Title\nReplace-Text \"Quantum Mysteries Unveiled\" \"Government and Public Policy: Unveiling the Interwoven Dynamics\"\n\n# 2) Byline (author name)\nReplace-Text \"Isaac Newton\" \"Dr. Henry Morgan\"\n\n# 3) Email line - originally split across runs (\"isaac\" + \".\" + \"newton@physics\" + \".\" + \"org\"),\n# which together read as \"isaac.newton@physics.org\"; replace the whole string.\nReplace-Text \"isaac.newton@physics.org\" \"hmorgan@moorehighschool.com\"\n\n# 4) First body paragraph - four sentences before the first line break.\nReplace-Text \"In the realm of physics, the exploration of quantum mechanics has yielded profound insights, challenging our understanding of the universe\" \"In the realm of human society, nothing is more pivotal than the intricate tapestry of government and public policy\"\nReplace-Text \" Quantum mechanics, initially conceived as a framework to explain the enigmatic behavior of subatomic particles, has revolutionized various disciplines, including chemistry, materials science, and even computer science\" \" These two forces, intertwined like threads in a multihued fabric, exert profound influences on the trajectory of our lives\"\nReplace-Text \" This essay delves into the captivating world of quantum physics, investigating its history, fundamental principles, and transformative applications\" \" Government shapes societal structure, crafting laws and policies that weave through every aspect of our existence\"\nInsert-After \"Government shapes societal structure, crafting laws and policies that weave through every aspect of our existence.\" \" Public policy, in turn, is a vivid reflection of the collective aspirations and values of society, a tapestry woven from the threads of governmental action.\"\n\n# Second block (after first double line-break)\nReplace-Text \"This conceptual revolution, sparked by Max Planck's introduction of energy quantization, shattered long-held assumptions about the continuity of energy and revealed the discreet nature of light and matter\" \"These interconnected entities paint a picture of complex interactions\"\nReplace-Text \" Albert Einstein's groundbreaking photoelectric effect experiment provided empirical validation for this departure from classical physics\" \" Government draws its legitimacy from the consent of the governed, while public policy becomes a testament to the collective will of the people\"\nReplace-Text \" Moreover, the advent of wave-particle duality, epitomized by the double-slit experiment, unveiled the paradoxical behavior of subatomic particles exhibiting both wave-like and particle-like properties\" \" As the government formulates policies that shape regulations, taxation, and public services, it lays the foundation for a just and equitable society\"\nInsert-After \"As the government formulates policies that shape regulations, taxation, and public services, it lays the foundation for a just and equitable society.\" \" These policies impact everything from economic growth and environmental protection to education and healthcare, weaving their way into the very fabric of our daily lives.\"\n\n# Third block (after second double line-break)\nReplace-Text \"The enigma of quantum entanglement, wherein particles exhibit a profound interconnectedness regardless of distance, has perplexed scientists and sparked profound debates about the nature of reality\" \"Furthermore, government and public policy are dynamic entities, constantly evolving in response to changing societal needs and global challenges\"\nReplace-Text \" This phenomenon, defying classical notions of locality, has profound implications, ranging from cryptography to quantum computing\" \" As the world grapples with issues such as climate change, resource depletion, and globalization, governments must adapt their policies to confront these pressing concerns\"\nReplace-Text \" It challenges our understanding of information transfer and raises fundamental questions regarding the relationship between consciousness and the physical world\" \" This delicate dance between government and public policy ensures that society navigates the ever-shifting tides of progress, preserving core values while embracing the imperatives of a changing world\"\n\n# 5) Summary paragraph\nReplace-Text \"Quantum mechanics has reshaped our understanding of the universe, providing a conceptual framework for explaining the enigmatic behavior of subatomic particles\" \"In this essay, we explored the intrinsic connection between government and public policy\"\nReplace-Text \" The quantization of energy, wave-particle duality, and quantum entanglement challenge classical notions of physics and open up new avenues for exploration across a wide range of fields\" \" We recognized the critical role government plays in shaping society through the policies it enacts\"\n# This sentence used to straddle a page-break run (\"... and \" + \"revolutionizing industries\"); now it\n# becomes one self-contained sentence, and the following sentence is replaced independently below.\nReplace-Text \" The transformative applications of quantum mechanics, from quantum computing to cryptography, hold immense promise for advancing technology and revolutionizing industries\" \" These policies, in turn, are influenced by societal values and aspirations, and form the backbone of a stable and just society\"\nReplace-Text \" As we delve deeper into the mysteries of the quantum realm, we unlock new frontiers of knowledge with the potential to shape our future in unimaginable ways\" \" Government and public policy work in tandem, evolving over time to address new Herausforderungen and societal shifts\"\nInsert-After \"Government and public policy work in tandem, evolving over time to address new Herausforderungen and societal shifts.\" \" By fostering this vital partnership, we ensure a government that is responsive to the people it serves and policies that reflect the collective will and best interests of society as a whole.\"\n\n# 6) Trailing empty paragraph added at the end of the document body.\n$endRng = $d.Content\n$endRng.Collapse(0)\n$endRng.InsertParagraphAfter()\n"}
